$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7530
$ws.Range("K3").Value = 7784
$ws.Range("K4").Value = 1638
$ws.Range("K6").Value = 8684
$ws.Range("K7").Value = 26187

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 477
$ws.Range("K3").Value = 515
$ws.Range("K6").Value = 573
$ws.Range("K7").Value = 1709

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 194
$ws.Range("K7").Value = 554

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 275
$ws.Range("K6").Value = 352
$ws.Range("K7").Value = 1101

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 149
$ws.Range("K7").Value = 434

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K6").Value = 262
$ws.Range("K7").Value = 866

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 200
$ws.Range("K7").Value = 613

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 117
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 439

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 227
$ws.Range("K4").Value = 92
$ws.Range("K7").Value = 789
$ws.Range("K8").Value = 1709
$ws.Range("K15").Value = 269
$ws.Range("K19").Value = 757
$ws.Range("K23").Value = 262
$ws.Range("K29").Value = 1442
$ws.Range("K31").Value = 308
$ws.Range("K33").Value = 1101
$ws.Range("K37").Value = 866
$ws.Range("K42").Value = 967
$ws.Range("K43").Value = 216
$ws.Range("K49").Value = 149
$ws.Range("K51").Value = 338
$ws.Range("K52").Value = 674
$ws.Range("K54").Value = 515
$ws.Range("K55").Value = 288
$ws.Range("K57").Value = 105
$ws.Range("K60").Value = 156
$ws.Range("K63").Value = 73
$ws.Range("K65").Value = 613
$ws.Range("K67").Value = 1021
$ws.Range("K69").Value = 60
$ws.Range("K72").Value = 125
$ws.Range("K75").Value = 83
$ws.Range("K77").Value = 169
$ws.Range("K78").Value = 327
$ws.Range("K79").Value = 643
$ws.Range("K83").Value = 554
$ws.Range("K85").Value = 1201
$ws.Range("K88").Value = 278
$ws.Range("K91").Value = 309
$ws.Range("K92").Value = 96
$ws.Range("K94").Value = 350
$ws.Range("K95").Value = 434
$ws.Range("K96").Value = 278
$ws.Range("K98").Value = 139
$ws.Range("K99").Value = 439
$ws.Range("K101").Value = 26187

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 308

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 367
$ws.Range("K6").Value = 291
$ws.Range("K7").Value = 1021

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 282
$ws.Range("K7").Value = 515

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 406
$ws.Range("K6").Value = 423
$ws.Range("K7").Value = 1442

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 257
$ws.Range("K7").Value = 757

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 365
$ws.Range("K7").Value = 967

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 99
$ws.Range("K4").Value = 31
$ws.Range("K7").Value = 327

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 82
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 79
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 309

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 203
$ws.Range("K7").Value = 643

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 259
$ws.Range("K3").Value = 248
$ws.Range("K6").Value = 221
$ws.Range("K7").Value = 789

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 73
$ws.Range("K6").Value = 163
$ws.Range("K7").Value = 350

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K3").Value = 19
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 73
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 92
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 338

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 22
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K6").Value = 296
$ws.Range("K7").Value = 1201

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 169

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 184
$ws.Range("K7").Value = 674

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 92
